# Updated RTL Description Doc
# changed i2si_en to rf_i2si_en
#
# The "i2si_en" signal row (Signal Name = i2si_en, Direction = in, Bits = 1,
# Comment = "i2s input is enabled") is removed from its original position
# (row 48 in the "i2si" block) and re-added, renamed to "rf_i2si_en", three
# rows further down -- immediately above the "rf_bist_en" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old "i2si_en" row entirely; rows below shift up.
$ws.Rows(48).Delete()

# 2) Insert a fresh row just above the row that now holds "rf_bist_en"
#    (this used to be row 52, now row 51 after the deletion above).
$ws.Rows(51).Insert()

# Copy formatting from the row above (so the new row matches the existing
# table styling exactly) without touching its values.
$ws.Range("A50:E50").Copy()
$ws.Range("A51:E51").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Populate the new row with the renamed signal.
$ws.Cells.Item(51, 2).Value = "rf_i2si_en"
$ws.Cells.Item(51, 3).Value = "in"
$ws.Cells.Item(51, 4).Value = 1
$ws.Cells.Item(51, 5).Value = "i2s input is enabled"

# 4) Reflect the user's resulting selection/scroll position.
$ws.Range("D52").Select()
